# Applies the diff:
#  - Insert a new row at row 4: 1004 | Test 4 | 1543 | PRJ-02
#    (this shifts the former rows 4 and 5 down to rows 5 and 6)
#  - Append a new row at the end (row 7): 1004 | Test 4 Dup1675 | (blank) | PRJ-02

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the existing row 4 ("1003, Test 2 Co, hf, PRJ-01"),
# pushing it (and the row after it) down by one.
$ws.Rows("4:4").Insert()

# Fill in the newly inserted row 4.
$ws.Cells.Item(4, 1).Value = 1004
$ws.Cells.Item(4, 2).Value = "Test 4"
$ws.Cells.Item(4, 3).Value = 1543
$ws.Cells.Item(4, 4).Value = "PRJ-02"

# Append a brand new row 7 at the bottom of the data.
$ws.Cells.Item(7, 1).Value = 1004
$ws.Cells.Item(7, 2).Value = "Test 4 Dup1675"
$ws.Cells.Item(7, 4).Value = "PRJ-02"
